# Update the "Totals" worksheet figures for the row dated 43739 (row 13).
# The running-month input figures changed; all dependent formula cells
# (D13, G13, H13, I13, B30, C30, D30, B31, C31) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Totals")

$ws.Range("B13").Value = 930
$ws.Range("C13").Value = 30315
$ws.Range("E13").Value = 18309
$ws.Range("F13").Value = 49582

# Move the sheet's selection from B30 to E40, matching the saved view state.
[void]$ws.Range("E40").Select()
